# Centra la finestra di incidenza a 7 giorni sull'ultimo giorno (trailing
# window [r-6, r]) invece che centrata sul giorno medio ([r-3, r+3]).
# Ricalcola le colonne C ("somma mobile 7gg.") e D ("somma mobile 7gg.
# per 100mila abitanti") per ogni riga di dati del foglio.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 184
$windowSize = 7
$population = 1861

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $winStart = $r - $windowSize + 1

    if ($winStart -lt $firstRow) {
        # Non ci sono ancora 7 giorni di dati disponibili: lascia vuoto.
        $ws.Cells.Item($r, 3).ClearContents()
        $ws.Cells.Item($r, 4).ClearContents()
    } else {
        $sum = 0
        for ($i = $winStart; $i -le $r; $i++) {
            $sum = $sum + $ws.Cells.Item($i, 2).Value2
        }
        $ws.Cells.Item($r, 3).Value = $sum
        $ws.Cells.Item($r, 4).Value = $sum * 100000.0 / $population
    }
}
